# Insert one new data row at row 686 (pushing the former rows 686..776 down to 687..777)
# and populate it with the new weekly observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 686, shifting existing rows 686..776 down to 687..777.
$ws.Rows.Item(686).Insert()

# Populate the newly inserted row 686 with the new record's data.
$ws.Range("A686").Value = 6
$ws.Range("B686").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C686").Value = "Metropolitana"
$ws.Range("D686").Value = 45124
$ws.Range("E686").Value = 13
$ws.Range("F686").Value = 100112044
$ws.Range("G686").Value = "Perejil"
$ws.Range("H686").Value = "Sin especificar"
$ws.Range("I686").Value = "Primera"
$ws.Range("J686").Value = 280
$ws.Range("K686").Value = 14000
$ws.Range("L686").Value = 15000
$ws.Range("M686").Value = 14464
$ws.Range("N686").Value = "`$/docena de atados"
$ws.Range("O686").Value = "Región Metropolitana"
$ws.Range("P686").Value = 4821
$ws.Range("Q686").Value = 3
$ws.Range("R686").Value = "Hortaliza"
